$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "56.672.15"
$ws.Range("E2").Value = "  +2.28%  "

$ws.Range("D3").Value = "2.323.95"
$ws.Range("E3").Value = "  +1.13%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "520.73"
$ws.Range("E5").Value = "  +2.69%  "

$ws.Range("D6").Value = "134.95"
$ws.Range("E6").Value = "  +3.73%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +1.27%  "

$ws.Range("D9").Value = "2.348.55"
$ws.Range("E9").Value = "  +1.10%  "

$ws.Range("E10").Value = "  +5.33%  "

$ws.Range("E11").Value = "  -0.78%  "

$ws.Range("E12").Value = "  +3.47%  "

$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("D14").Value = "23.92"
$ws.Range("E14").Value = "  -0.14%  "

$ws.Range("D15").Value = "2.742.52"
$ws.Range("E15").Value = "  +1.31%  "

$ws.Range("D16").Value = "56.732.91"
$ws.Range("E16").Value = "  +2.73%  "

$ws.Range("E17").Value = "  +2.08%  "

$ws.Range("D18").Value = "2.345.63"
$ws.Range("E18").Value = "  -7.63%  "

$ws.Range("E19").Value = "  -2.88%  "

$ws.Range("D20").Value = "4.23"
$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("D21").Value = "323.06"
$ws.Range("E21").Value = "  +3.85%  "

$ws.Range("D22").Value = "6.61"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").Value = "60.57"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").Value = "0.165"
$ws.Range("E25").Value = "  +8.47%  "

$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").Value = "7.93"
$ws.Range("E27").Value = "  +5.23%  "

$ws.Range("E28").Value = "  +12.61%  "

$ws.Range("E29").Value = "  +5.18%  "

$ws.Range("D30").Value = "169.18"
$ws.Range("E30").Value = "  -2.24%  "

$ws.Range("E31").Value = "  +5.09%  "

$ws.Range("D32").Value = "6.19"
$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").Value = "0.992"
$ws.Range("E35").Value = "  -0.41%  "

$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("E38").Value = "  +3.53%  "

$ws.Range("E39").Value = "  +7.16%  "

$ws.Range("D40").Value = "37.90"
$ws.Range("E40").Value = "  +2.91%  "

$ws.Range("E41").Value = "  +0.34%  "

$ws.Range("D42").Value = "3.60"
$ws.Range("E42").Value = "  +4.33%  "

$ws.Range("D43").Value = "138.46"
$ws.Range("E43").Value = "  +2.37%  "

$ws.Range("D44").Value = "5.28"
$ws.Range("E44").Value = "  +6.66%  "

$ws.Range("D45").Value = "277.26"
$ws.Range("E45").Value = "  +5.95%  "

$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("D47").Value = "0.0507"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").Value = "0.563"
$ws.Range("E48").Value = "  +1.70%  "

$ws.Range("E49").Value = "  +3.24%  "

$ws.Range("D50").Value = "17.91"
$ws.Range("E50").Value = "  +7.41%  "

$ws.Range("E51").Value = "  +0.34%  "
